# update synapse table all except Prostate
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The CRC v2.0-public rows (72:79) all had release_date values that were
# placeholder/incorrect ("2022-12".."2022-18"); they are corrected here to
# the actual CRC v2.0-public release date, matching row 71 ("2022-10").
$crcReleaseDates = New-Object 'object[,]' 8,1
for ($i = 0; $i -lt 8; $i++) { $crcReleaseDates[$i,0] = "2022-10" }
$ws.Range("E72:E79").Value = $crcReleaseDates

# Append new cohort rows: PANC, Prostate, BLADDER, and an updated BrCa
# (v1.2-consortium) block with a ca_radtx_dataset row and release dates.
$newRows = New-Object 'object[,]' 36,5

$data = @(
  @("PANC", "v1.1-consortium", "cancer_level_dataset_index", "syn47283323", "2022-02"),
  @("PANC", "v1.1-consortium", "cancer_level_dataset_non_index", "syn47283351", "2022-02"),
  @("PANC", "v1.1-consortium", "cancer_panel_test_level_dataset", "syn47283590", "2022-02"),
  @("PANC", "v1.1-consortium", "imaging_level_dataset", "syn47283470", "2022-02"),
  @("PANC", "v1.1-consortium", "med_onc_note_level_dataset", "syn47283544", "2022-02"),
  @("PANC", "v1.1-consortium", "pathology_report_level_dataset", "syn47283508", "2022-02"),
  @("PANC", "v1.1-consortium", "patient_level_dataset", "syn47283390", "2022-02"),
  @("PANC", "v1.1-consortium", "regimen_cancer_level_dataset", "syn47283430", "2022-02"),
  @("PANC", "v1.1-consortium", "tm_level_dataset", "syn47283628", "2022-02"),
  @("Prostate", "v1.1-consortium", "cancer_level_dataset_index", "", ""),
  @("Prostate", "v1.1-consortium", "cancer_level_dataset_non_index", "", ""),
  @("Prostate", "v1.1-consortium", "cancer_panel_test_level_dataset", "", ""),
  @("Prostate", "v1.1-consortium", "imaging_level_dataset", "", ""),
  @("Prostate", "v1.1-consortium", "med_onc_note_level_dataset", "", ""),
  @("Prostate", "v1.1-consortium", "pathology_report_level_dataset", "", ""),
  @("Prostate", "v1.1-consortium", "patient_level_dataset", "", ""),
  @("Prostate", "v1.1-consortium", "regimen_cancer_level_dataset", "", ""),
  @("Prostate", "v1.1-consortium", "tm_level_dataset", "", ""),
  @("BLADDER", "v1.1-consortium", "ca_radtx_dataset", "syn44420748", "2022-11"),
  @("BLADDER", "v1.1-consortium", "cancer_level_dataset_index", "syn44420702", "2022-11"),
  @("BLADDER", "v1.1-consortium", "cancer_level_dataset_non_index", "syn44420708", "2022-11"),
  @("BLADDER", "v1.1-consortium", "cancer_panel_test_level_dataset", "syn44420744", "2022-11"),
  @("BLADDER", "v1.1-consortium", "imaging_level_dataset", "syn44420731", "2022-11"),
  @("BLADDER", "v1.1-consortium", "med_onc_note_level_dataset", "syn44420739", "2022-11"),
  @("BLADDER", "v1.1-consortium", "pathology_report_level_dataset", "syn44420737", "2022-11"),
  @("BLADDER", "v1.1-consortium", "patient_level_dataset", "syn44420719", "2022-11"),
  @("BLADDER", "v1.1-consortium", "regimen_cancer_level_dataset", "syn44420726", "2022-11"),
  @("BrCa", "v1.2-consortium", "cancer_level_dataset_index", "syn43172806", "2022-10"),
  @("BrCa", "v1.2-consortium", "cancer_level_dataset_non_index", "syn43172815", "2022-10"),
  @("BrCa", "v1.2-consortium", "cancer_panel_test_level_dataset", "syn43172901", "2022-10"),
  @("BrCa", "v1.2-consortium", "imaging_level_dataset", "syn43172865", "2022-10"),
  @("BrCa", "v1.2-consortium", "med_onc_note_level_dataset", "syn43172895", "2022-10"),
  @("BrCa", "v1.2-consortium", "pathology_report_level_dataset", "syn43172879", "2022-10"),
  @("BrCa", "v1.2-consortium", "patient_level_dataset", "syn43172821", "2022-10"),
  @("BrCa", "v1.2-consortium", "regimen_cancer_level_dataset", "syn43172821", "2022-10"),
  @("BrCa", "v1.2-consortium", "tm_level_dataset", "syn43172910", "2022-10")
)

for ($i = 0; $i -lt $data.Count; $i++) {
  for ($j = 0; $j -lt 5; $j++) {
    $newRows[$i, $j] = $data[$i][$j]
  }
}

$ws.Range("A80:E115").Value = $newRows

# Leave the view scrolled/selected near the newly-added rows, as in the
# authored workbook.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 73
$ws.Range("E89:E97").Select()
